$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2473.8333
$ws.Range("I38").Value = 873.6
$ws.Range("J38").Value = 3616.8572
$ws.Range("K38").Value = 2620.8
$ws.Range("L38").Value = 10850.5716
$ws.Range("M38").Value = -2248.8
$ws.Range("N38").Value = -11594.5716

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6538287.5
$ws.Range("J51").Value = 7938563.5
$ws.Range("L51").Value = 7938563.5
$ws.Range("N51").Value = -7939531.5

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 27944.445
$ws.Range("J87").Value = 27944.445
$ws.Range("L87").Value = 27944.445
$ws.Range("N87").Value = -30440.445

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 27944.445
$ws.Range("J90").Value = 27944.445
$ws.Range("L90").Value = 83833.33499999999
$ws.Range("N90").Value = -96313.33499999999

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1378.75
$ws.Range("J112").Value = 1636.1818
$ws.Range("L112").Value = 4908.5454
$ws.Range("N112").Value = -7124.5454

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 942.5925999999999
$ws.Range("I129").Value = 380
$ws.Range("J129").Value = 1103.3334
$ws.Range("K129").Value = 1140
$ws.Range("L129").Value = 3310.0002
$ws.Range("M129").Value = 3860
$ws.Range("N129").Value = -13310.0002

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2629.3582
$ws.Range("I138").Value = 1405.4667
$ws.Range("J138").Value = 2982.4038
$ws.Range("K138").Value = 4216.4001
$ws.Range("L138").Value = 8947.2114
$ws.Range("M138").Value = 923.5999000000002
$ws.Range("N138").Value = -19227.2114

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4542.3335
$ws.Range("I141").Value = 2344.5833
$ws.Range("J141").Value = 13333.333
$ws.Range("K141").Value = 7033.749899999999
$ws.Range("L141").Value = 39999.999
$ws.Range("M141").Value = -1853.749899999999
$ws.Range("N141").Value = -50359.999

# ARM row 62
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# ARM row 65
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# ARM row 81
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H81").Value = 41400
$ws.Range("J81").Value = 41400
$ws.Range("L81").Value = 41400
$ws.Range("N81").Value = -43396

# ARM row 84
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H84").Value = 41400
$ws.Range("J84").Value = 41400
$ws.Range("L84").Value = 124200
$ws.Range("N84").Value = -134184

# ARM row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 28612.5
$ws.Range("J104").Value = 28612.5
$ws.Range("L104").Value = 28612.5
$ws.Range("N104").Value = -35600.5

# BSM row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 664.9211
$ws.Range("I5").Value = 405.75
$ws.Range("J5").Value = 1109.2142
$ws.Range("K5").Value = 1217.25
$ws.Range("L5").Value = 3327.6426
$ws.Range("M5").Value = -1105.25
$ws.Range("N5").Value = -3551.6426

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 94.25
$ws.Range("I23").Value = 76.666664
$ws.Range("J23").Value = 104.8
$ws.Range("K23").Value = 229.999992
$ws.Range("L23").Value = 314.4
$ws.Range("M23").Value = 5.000008000000008
$ws.Range("N23").Value = -784.4

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 827.05884
$ws.Range("I122").Value = 360.18182
$ws.Range("K122").Value = 3241.63638
$ws.Range("M122").Value = -791.6363799999999

# CUL row 130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 871.5
$ws.Range("I130").Value = 871.5
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 2614.5
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 2405.5
$ws.Range("N130").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1905772.2
$ws.Range("I131").Value = 16667154
$ws.Range("J131").Value = 1077.8387
$ws.Range("K131").Value = 50001462
$ws.Range("L131").Value = 3233.5161
$ws.Range("M131").Value = -49996422
$ws.Range("N131").Value = -13313.5161

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 664.9211
$ws.Range("I135").Value = 405.75
$ws.Range("J135").Value = 1109.2142
$ws.Range("K135").Value = 3651.75
$ws.Range("L135").Value = 9982.927799999999
$ws.Range("M135").Value = -1116.75
$ws.Range("N135").Value = -15052.9278

# LTW row 69
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 22500
$ws.Range("J69").Value = 22500
$ws.Range("L69").Value = 22500
$ws.Range("N69").Value = -24122

# LTW row 72
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 22500
$ws.Range("J72").Value = 22500
$ws.Range("L72").Value = 67500
$ws.Range("N72").Value = -75612

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 33000
$ws.Range("J140").Value = 33000
$ws.Range("L140").Value = 33000
$ws.Range("N140").Value = -43360

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 119536.88
$ws.Range("I62").Value = 4134.737
$ws.Range("J62").Value = 432771.28
$ws.Range("K62").Value = 4134.737
$ws.Range("L62").Value = 432771.28
$ws.Range("M62").Value = -3510.737
$ws.Range("N62").Value = -434019.28

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 119536.88
$ws.Range("I65").Value = 4134.737
$ws.Range("J65").Value = 432771.28
$ws.Range("K65").Value = 20673.685
$ws.Range("L65").Value = 2163856.4
$ws.Range("M65").Value = -17553.685
$ws.Range("N65").Value = -2170096.4

# WVR row 68
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 50271
$ws.Range("J68").Value = 50271
$ws.Range("L68").Value = 50271
$ws.Range("N68").Value = -51893

# WVR row 71
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 50271
$ws.Range("J71").Value = 50271
$ws.Range("L71").Value = 150813
$ws.Range("N71").Value = -158925

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 46570
$ws.Range("I100").Value = 54655.453
$ws.Range("J100").Value = 2100
$ws.Range("K100").Value = 109310.906
$ws.Range("L100").Value = 4200
$ws.Range("M100").Value = -108769.906
$ws.Range("N100").Value = -5282

# WVR row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 47057.145
$ws.Range("J138").Value = 47057.145
$ws.Range("L138").Value = 47057.145
$ws.Range("N138").Value = -57337.145
